# Web Deployment Ver 1.0
# Updated Shiny App for web deployment and test datasets with realistic data.
# Updated calculation to correct formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old placeholder test data (Depth/Area) with the new,
# realistic test dataset produced by the corrected formula.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2790

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 6000

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 12000

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 20000

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 30492

# Leave the selection where the author left it when they saved the file.
$ws.Range("D7").Select()
